$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: ledipasvir / LDV / NS5A inhibitors / Gilead Sciences / GS-5885
# Row 9: daclatasvir / DCV / NS5A inhibitors / Bristol-Myers Squibb  / BMS-790052
$ws.Range("A8").Value = "ledipasvir"
$ws.Range("A9").Value = "daclatasvir"
$ws.Range("B9").Value = "DCV"
$ws.Range("B8").Value = "LDV"
$ws.Range("E8").Value = "GS-5885"
$ws.Range("D9").Value = "Bristol-Myers Squibb "
$ws.Range("E9").Value = "BMS-790052"
$ws.Range("C8").Value = "NS5A inhibitors"
$ws.Range("C9").Value = "NS5A inhibitors"
$ws.Range("D8").Value = "Gilead Sciences"

$ws.Range("A1:XFD1048576").Select()
